$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "Sayfa1" to "Data"
$ws.Name = "Data"

# --- Update cell values (objective 3 scaling) ---

# Rows 3-8 (A=1): D,E,F were 8,8,8 -> now 4,2,1
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("D5").Value = 4
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 1
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 1
$ws.Range("D7").Value = 4
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 1
$ws.Range("D8").Value = 4
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 1

# Rows 9-14 (A=2): D,E,F were 15,15,15 -> now 6,3,2
$ws.Range("D9").Value = 6
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 2
$ws.Range("D10").Value = 6
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 2
$ws.Range("D11").Value = 6
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 2
$ws.Range("D12").Value = 6
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 2
$ws.Range("D13").Value = 6
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 2
$ws.Range("D14").Value = 6
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 2

# Rows 15-20 (A=3): G,H,I were 15,15,15 -> now 10,50,50
$ws.Range("G15").Value = 10
$ws.Range("H15").Value = 50
$ws.Range("I15").Value = 50
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 50
$ws.Range("I16").Value = 50
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = 50
$ws.Range("I17").Value = 50
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 50
$ws.Range("G19").Value = 10
$ws.Range("H19").Value = 50
$ws.Range("I19").Value = 50
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 50
$ws.Range("I20").Value = 50

# Rows 21-26 (A=4): G,H,I were 8,8,8 -> now 7,20,12
$ws.Range("G21").Value = 7
$ws.Range("H21").Value = 20
$ws.Range("I21").Value = 12
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = 20
$ws.Range("I22").Value = 12
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 20
$ws.Range("I23").Value = 12
$ws.Range("G24").Value = 7
$ws.Range("H24").Value = 20
$ws.Range("I24").Value = 12
$ws.Range("G25").Value = 7
$ws.Range("H25").Value = 20
$ws.Range("I25").Value = 12
$ws.Range("G26").Value = 7
$ws.Range("H26").Value = 20
$ws.Range("I26").Value = 12

# --- View state: zoom to 125% and move selection to I18 ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 125
$ws.Range("I18").Select()
